$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 4.429675500412797 }
    3 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732;  G = 3.781711156805759 }
    4 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 3.645393585217082 }
    5 = @{ B = 0.6753301551942219; C = 0.3127903958511391; D = 0.8054896365839992;  E = 0.496779210170732;  G = 2.290389397800092 }
    6 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 6.201049113329182 }
    7 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.1575252929769615;  E = 0.496779210170732;  G = 5.553084769722144 }
    8 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 6.201049113329182 }
    9 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992;  E = 0.496779210170732;  G = 6.201049113329182 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
